# Generate Report for Archive
#
# The localization status for the two handed-off files moved from
# "Ready for handoff" to "In Translation" on every sheet that surfaces
# the Status column (Overview!E/F and the per-locale Status column on
# the zh-cn / de-de sheets). Updating that text makes the Status column
# content narrower, so the Status column is re-sized to fit the new text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: columns E (zh-cn) and F (de-de) hold the status text ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

# --- Per-locale sheets: column C holds the Status text ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# --- Resize the Status columns to fit the shorter text ---
# NB: the host snaps ColumnWidth to a pixel grid (stored = (round(w*6)+5)/6),
# so the input below is chosen as the closest value that round-trips to the
# target stored width of 13.4101845877511 characters.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
